# Advances on consolidated report.
# The error-report block that used to start at "countryadmin1" (old col AE) grew
# a new leading "missingcountry" check column, and the old lead-in columns
# (countryadmin1, Admin1and2, sectorindicator, Indicatortype -> old AE:AH)
# were dropped. That shifts everything from the old "missingcountry" column
# (AI) onward 4 columns to the left, and the final summary column's header
# text changes from "ERROR" to "Review".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four now-obsolete lead-in columns (old AE:AH); Excel shifts
# everything to the right of them left by 4 columns automatically.
$ws.Range("AE1:AH19").EntireColumn.Delete()

# The last column's header label changes from "ERROR" to "Review".
$ws.Range("AX1").Value = "Review"
